# Update "想去人数" (column F) counts across the four sheets of the
# "上海-漫展信息" workbook to reflect newly generated output.
#
# Sheet order in the workbook:
#   1 -> 展览      (Exhibitions)
#   2 -> 演出      (Performances)
#   3 -> 本地生活  (Local life)
#   4 -> 全部类型  (All types - aggregate of the other sheets)

$wb = $excel.ActiveWorkbook

function Set-F {
    param($SheetIndex, $Row, $Value)
    $ws = $wb.Worksheets.Item($SheetIndex)
    $ws.Cells.Item($Row, 6).Value = $Value
}

# 展览 (Sheet 1)
Set-F 1 9  570
Set-F 1 13 1366
Set-F 1 17 117
Set-F 1 25 1175
Set-F 1 27 24
Set-F 1 35 876
Set-F 1 45 101

# 演出 (Sheet 2)
Set-F 2 13 2467
Set-F 2 14 1167
Set-F 2 26 46445
Set-F 2 45 109

# 本地生活 (Sheet 3)
Set-F 3 6  2668
Set-F 3 7  4446
Set-F 3 12 387
Set-F 3 13 120
Set-F 3 16 232

# 全部类型 (Sheet 4)
Set-F 4 5  4446
Set-F 4 9  120
Set-F 4 10 120
Set-F 4 19 570
Set-F 4 22 2467
Set-F 4 23 1167
Set-F 4 24 1366
Set-F 4 27 117
Set-F 4 34 1175
Set-F 4 40 876
Set-F 4 50 101
